# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the rows whose
# handoff was just (re-)generated: rows 4, 6, 7, 8, 9, 10 on both the
# "zh-cn" and "de-de" sheets. Rows 2/3 (already handed back) and row 5
# (in translation) keep their existing handoff timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value = "2016-02-18 08:42:58"
$zhcn.Range("D6").Value = "2016-02-18 08:42:58"
$zhcn.Range("D7").Value = "2016-02-18 08:42:58"
$zhcn.Range("D8").Value = "2016-02-18 08:42:58"
$zhcn.Range("D9").Value = "2016-02-18 08:42:58"
$zhcn.Range("D10").Value = "2016-02-18 08:42:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value = "2016-02-18 08:43:11"
$dede.Range("D6").Value = "2016-02-18 08:43:11"
$dede.Range("D7").Value = "2016-02-18 08:43:11"
$dede.Range("D8").Value = "2016-02-18 08:43:11"
$dede.Range("D9").Value = "2016-02-18 08:43:11"
$dede.Range("D10").Value = "2016-02-18 08:43:11"
